# Add a new data row (row 3) to the active worksheet, mirroring the
# structure of the existing header (row 1) / data (row 2) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cells -------------------------------------------------
$ws.Range("A3").Value = 131143915
$ws.Range("B3").Value = 58043
$ws.Range("E3").Value = 103021
$ws.Range("Q3").Value = 562303
$ws.Range("R3").Value = 6917048
$ws.Range("S3").Value = 10

# --- Text cells ------------------------------------------------------
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"

# I3 holds the text "4" (not a numeric value) in the source data, so
# force the cell to be treated as text before assigning it.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "4"

$ws.Range("M3").Value = "födosökande"
$ws.Range("N3").Value = "observerad"
$ws.Range("P3").Value = "Storbackmyran, Mpd"
$ws.Range("T3").Value = "Västernorrland"
$ws.Range("U3").Value = "Ånge"
$ws.Range("V3").Value = "Medelpad"
$ws.Range("W3").Value = "Torp"

# Y3/Z3/AA3/AB3 store dates/times as plain text in the source data, so
# force text formatting before assigning to avoid Excel auto-converting
# them to date/time serial values.
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2026-01-03"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "14:10"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2026-01-03"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "14:13"

$ws.Range("AW3").Value = "Markus Borja"
$ws.Range("AX3").Value = "Markus Borja"

# --- Boolean cells ---------------------------------------------------
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

Write-Output "Row 3 added"
